$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Users" column (D) to hold "State".
$ws.Columns.Item(4).Insert()

# Header for the new column.
$ws.Range("D1").Value = "State"

# Users column text (same for all three rows per the diff).
$usersText = "root`ndhcpcd`ntss`npollinate`nubuntu`n"

# Row 2
$ws.Range("A2").Value = "i-056820ee844baf41e"
$ws.Range("B2").Value = ""
$ws.Range("D2").Value = "running"
$ws.Range("E2").Value = $usersText

# Row 3
$ws.Range("A3").Value = "i-0c4995636b4a3f8ad"
$ws.Range("B3").Value = ""
$ws.Range("D3").Value = "running"
$ws.Range("E3").Value = $usersText

# Row 4
$ws.Range("A4").Value = "i-046e18bd3080b03ce"
$ws.Range("B4").Value = ""
$ws.Range("D4").Value = "running"
$ws.Range("E4").Value = $usersText

# The multi-line Users text makes the host auto-expand row height; restore
# the default (un-customized) row heights so the sheet metadata matches.
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
